# Update odds values on Sheet1 to reflect the latest FlashScore snapshot
# for Jogos_da_Semana_FlashScore_2025-04-15.xlsx

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 4
$ws.Range("G4").Value = 2.63
$ws.Range("H4").Value = 3.1
$ws.Range("I4").Value = 2.8
$ws.Range("AB4").Value = 17

# Row 5
$ws.Range("N5").Value = 1.53
$ws.Range("O5").Value = 2.4

# Row 6
$ws.Range("P6").Value = 1.33
$ws.Range("Q6").Value = 3.25
$ws.Range("U6").Value = 23
$ws.Range("AE6").Value = 8.5

# Row 7
$ws.Range("N7").Value = 1.73
$ws.Range("O7").Value = 2.08

# Row 9
$ws.Range("J9").Value = 1.1
$ws.Range("K9").Value = 7

# Row 13
$ws.Range("J13").Value = 1.07
$ws.Range("K13").Value = 9

# Row 14
$ws.Range("N14").Value = 1.95
$ws.Range("O14").Value = 1.85

# Row 15
$ws.Range("H15").Value = 3.3
$ws.Range("I15").Value = 2.25
$ws.Range("J15").Value = 1.06
$ws.Range("K15").Value = 10
$ws.Range("U15").Value = 15
$ws.Range("X15").Value = 23
$ws.Range("AG15").Value = 9.5
$ws.Range("AI15").Value = 19

# Row 16
$ws.Range("N16").Value = 1.43
$ws.Range("R16").Value = 2.65
$ws.Range("S16").Value = 1.43

# Row 17
$ws.Range("H17").Value = 3.6
$ws.Range("I17").Value = 3.35
$ws.Range("O17").Value = 2
$ws.Range("T17").Value = 7.5
$ws.Range("X17").Value = 11.5
$ws.Range("Y17").Value = 17.5
$ws.Range("Z17").Value = 12.5
$ws.Range("AA17").Value = 6.3
$ws.Range("AB17").Value = 11
$ws.Range("AC17").Value = 37
$ws.Range("AD17").Value = 200
$ws.Range("AE17").Value = 10.25
$ws.Range("AF17").Value = 16
$ws.Range("AG17").Value = 10
$ws.Range("AJ17").Value = 24

# Row 20
$ws.Range("O20").Value = 1.59

# Row 21
$ws.Range("G21").Value = 3.3
$ws.Range("I21").Value = 2.3
$ws.Range("O21").Value = 1.47
$ws.Range("X21").Value = 29
$ws.Range("AF21").Value = 10
